$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-08-16 10:45:47"
$wsZhCn.Range("H4").Value = "2016-08-16 10:45:42"
$wsZhCn.Range("K4").Value = "2016-08-16 10:46:02"
$wsDeDe.Range("H4").Value = "2016-08-16 10:45:47"
$wsDeDe.Range("K4").Value = "2016-08-16 10:46:14"
